{"js": "// The Introduction Letter used the placeholder \"YYY\" in two places where\n// the apprentice's full (mocking) title should have appeared instead, e.g.\n// \"...have sent you, YYY to the Looping Infinite Dungeon...\" and\n// \"As you, YYY, cannot hope to match the brilliance of myself, ...\".\n// Replace every \"YYY\" placeholder with the apprentice title phrase that is\n// already used elsewhere in the letter.\nconst apprenticeTitle =\n  \"my annoying, babbling, clumsy, dim-witted, elephantine, fool-hardy,  worthless apprentice (whoever or whatever you are)\";\n\nconst body = context.document.body;\nconst results = body.search(\"YYY\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(apprenticeTitle, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The Introduction Letter used the placeholder \"YYY\" in two places where\n# the apprentice's full (mocking) title should have appeared instead, e.g.\n# \"...have sent you, YYY to the Looping Infinite Dungeon...\" and\n# \"As you, YYY, cannot hope to match the brilliance of myself, ...\".\n# Replace every \"YYY\" placeholder with the apprentice title phrase that is\n# already used elsewhere in the letter.\n\n$apprenticeTitle = \"my annoying, babbling, clumsy, dim-witted, elephantine, fool-hardy,  worthless apprentice (whoever or whatever you are)\"\n\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"YYY\"\n$find.Replacement.Text = $apprenticeTitle\n$find.Execute([ref]\"YYY\", $false, $true, $false, $false, $false, $true, 1, $false, $apprenticeTitle, 2)\n"}
